{"js": "// Replace the date line and each two-digit multiplication answer cell\n// with its updated value, matching the target OOXML diff exactly.\nconst replacements = [\n  [\"2024-02-26 Monday\", \"2024-02-27 Tuesday\"],\n  [\"37\u00d788=3256\", \"66\u00d788=5808\"],\n  [\"49\u00d723=1127\", \"40\u00d721=840\"],\n  [\"20\u00d739=780\", \"99\u00d724=2376\"],\n  [\"81\u00d774=5994\", \"37\u00d737=1369\"],\n  [\"37\u00d773=2701\", \"12\u00d762=744\"],\n  [\"33\u00d796=3168\", \"88\u00d778=6864\"],\n  [\"43\u00d763=2709\", \"50\u00d768=3400\"],\n  [\"55\u00d723=1265\", \"73\u00d747=3431\"],\n  [\"28\u00d713=364\", \"94\u00d749=4606\"],\n  [\"42\u00d761=2562\", \"37\u00d794=3478\"],\n  [\"81\u00d746=3726\", \"15\u00d715=225\"],\n  [\"71\u00d744=3124\", \"33\u00d756=1848\"],\n  [\"82\u00d765=5330\", \"82\u00d783=6806\"],\n  [\"64\u00d783=5312\", \"48\u00d722=1056\"],\n  [\"20\u00d790=1800\", \"43\u00d789=3827\"],\n  [\"15\u00d735=525\", \"68\u00d786=5848\"],\n  [\"76\u00d740=3040\", \"81\u00d752=4212\"],\n  [\"42\u00d731=1302\", \"47\u00d768=3196\"],\n  [\"27\u00d726=702\", \"39\u00d771=2769\"],\n  [\"86\u00d737=3182\", \"69\u00d790=6210\"],\n  [\"77\u00d728=2156\", \"87\u00d727=2349\"],\n  [\"25\u00d736=900\", \"54\u00d783=4482\"],\n  [\"36\u00d729=1044\", \"58\u00d794=5452\"],\n  [\"12\u00d793=1116\", \"41\u00d742=1722\"],\n  [\"78\u00d742=3276\", \"66\u00d740=2640\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const item of found.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and each two-digit multiplication answer cell\n# to its new value via Find/Replace, matching the target OOXML diff.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-02-26 Monday\", \"2024-02-27 Tuesday\"),\n    @(\"37\u00d788=3256\", \"66\u00d788=5808\"),\n    @(\"49\u00d723=1127\", \"40\u00d721=840\"),\n    @(\"20\u00d739=780\", \"99\u00d724=2376\"),\n    @(\"81\u00d774=5994\", \"37\u00d737=1369\"),\n    @(\"37\u00d773=2701\", \"12\u00d762=744\"),\n    @(\"33\u00d796=3168\", \"88\u00d778=6864\"),\n    @(\"43\u00d763=2709\", \"50\u00d768=3400\"),\n    @(\"55\u00d723=1265\", \"73\u00d747=3431\"),\n    @(\"28\u00d713=364\", \"94\u00d749=4606\"),\n    @(\"42\u00d761=2562\", \"37\u00d794=3478\"),\n    @(\"81\u00d746=3726\", \"15\u00d715=225\"),\n    @(\"71\u00d744=3124\", \"33\u00d756=1848\"),\n    @(\"82\u00d765=5330\", \"82\u00d783=6806\"),\n    @(\"64\u00d783=5312\", \"48\u00d722=1056\"),\n    @(\"20\u00d790=1800\", \"43\u00d789=3827\"),\n    @(\"15\u00d735=525\", \"68\u00d786=5848\"),\n    @(\"76\u00d740=3040\", \"81\u00d752=4212\"),\n    @(\"42\u00d731=1302\", \"47\u00d768=3196\"),\n    @(\"27\u00d726=702\", \"39\u00d771=2769\"),\n    @(\"86\u00d737=3182\", \"69\u00d790=6210\"),\n    @(\"77\u00d728=2156\", \"87\u00d727=2349\"),\n    @(\"25\u00d736=900\", \"54\u00d783=4482\"),\n    @(\"36\u00d729=1044\", \"58\u00d794=5452\"),\n    @(\"12\u00d793=1116\", \"41\u00d742=1722\"),\n    @(\"78\u00d742=3276\", \"66\u00d740=2640\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
